# Path to Graduation - add student header info and extend the planner
# with a 2024 (Fall/Spring/Summer) block, restructuring the existing rows
# into three aligned 9-row blocks (2022, 2023, 2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for cell content (keeps existing formatting,
# e.g. the bold/centered title style on row 2, intact).
$ws.Cells.ClearContents()

# --- Header (new row 1): student name + id ---
$ws.Range("C1").Value = "Philip"
# The ID looks numeric, so force it to be stored as text the same way the
# original author's "394728739813" entry is (no leftover custom style).
$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "394728739813"
$ws.Range("E1").ClearFormats()

# --- Title row (row 2, already styled/merged A2:F2) ---
$ws.Range("A2").Value = "Path To Graduation"
# Re-writing the cell makes the engine recompute an explicit row height
# (because of the large title font); auto-fit restores the default
# (non-custom) row height like the original file had.
$ws.Rows("2:2").AutoFit()

# --- 2022 block (rows 3-11) ---
$ws.Range("A3").Value = "Fall 2022"
$ws.Range("B3").Value = "Credits"
$ws.Range("C3").Value = "Spring 2022"
$ws.Range("D3").Value = "Credits"
$ws.Range("E3").Value = "Summer 2022"
$ws.Range("F3").Value = "Credits"

$ws.Range("A4").Value = "POLS 1101"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "CPSC 3165"
$ws.Range("D4").Value = 3

$ws.Range("A5").Value = "DSCI 3111"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "CPSC 4135"
$ws.Range("D5").Value = 3

$ws.Range("A6").Value = "CPSC 3121"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CPSC 4148"
$ws.Range("D6").Value = 3

$ws.Range("A7").Value = "CPSC 4000"
$ws.Range("B7").Value = 0

$ws.Range("A11").Value = "Total"
$ws.Range("B11").Formula = "=SUM(B4:B10)"
$ws.Range("C11").Value = "Total"
$ws.Range("D11").Formula = "=SUM(D4:D10)"
$ws.Range("E11").Value = "Total"
$ws.Range("F11").Formula = "=SUM(F4:F10)"

# --- 2023 block (rows 12-20) ---
$ws.Range("A12").Value = "Fall 2023"
$ws.Range("B12").Value = "Credits"
$ws.Range("C12").Value = "Spring 2023"
$ws.Range("D12").Value = "Credits"
$ws.Range("E12").Value = "Summer 2023"
$ws.Range("F12").Value = "Credits"

$ws.Range("A13").Value = "CPSC 4155"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CPSC 4176"
$ws.Range("D13").Value = 3

$ws.Range("A14").Value = "CPSC 4157"
$ws.Range("B14").Value = 3

$ws.Range("A15").Value = "CPSC 4175"
$ws.Range("B15").Value = 3

$ws.Range("A20").Value = "Total"
$ws.Range("B20").Formula = "=SUM(B13:B19)"
$ws.Range("C20").Value = "Total"
$ws.Range("D20").Formula = "=SUM(D13:D19)"
$ws.Range("E20").Value = "Total"
$ws.Range("F20").Formula = "=SUM(F13:F19)"

# --- 2024 block (rows 21-29) ---
$ws.Range("A21").Value = "Fall 2024"
$ws.Range("B21").Value = "Credits"
$ws.Range("C21").Value = "Spring 2024"
$ws.Range("D21").Value = "Credits"
$ws.Range("E21").Value = "Summer 2024"
$ws.Range("F21").Value = "Credits"

$ws.Range("A29").Value = "Total"
$ws.Range("B29").Formula = "=SUM(B22:B28)"
$ws.Range("C29").Value = "Total"
$ws.Range("D29").Formula = "=SUM(D22:D28)"
$ws.Range("E29").Value = "Total"
$ws.Range("F29").Formula = "=SUM(F22:F28)"
